$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Ryk"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.43424333333333
$ws.Range("H2").Value = 31.30273
$ws.Range("I2").Value = 0.9711091978791583
$ws.Range("J2").Value = 0.9711091978791584
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 5.611633666666667
$ws.Range("N2").Value = 16.834901
$ws.Range("O2").Value = 0.1044932796553548
$ws.Range("P2").Value = 0.1044932796553548
$ws.Range("Q2").Value = 58.55315117552556
$ws.Range("R2").Value = 526.9783605797301
$ws.Range("S2").Value = 0.1014743849898741
$ws.Range("T2").Value = 0.1014743849898741

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Ryk"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.43424333333333
$ws.Range("H3").Value = 31.30273
$ws.Range("I3").Value = 0.9711091978791583
$ws.Range("J3").Value = 0.9711091978791584
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 23.41023966666667
$ws.Range("N3").Value = 70.230719
$ws.Range("O3").Value = 0.435918106133421
$ws.Range("P3").Value = 0.435918106133421
$ws.Range("Q3").Value = 244.2681371736522
$ws.Range("R3").Value = 2198.41323456287
$ws.Range("S3").Value = 0.4233240823882283
$ws.Range("T3").Value = 0.4233240823882283

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Ryk"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.43424333333333
$ws.Range("H4").Value = 31.30273
$ws.Range("I4").Value = 0.9711091978791583
$ws.Range("J4").Value = 0.9711091978791584
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1135936666666667
$ws.Range("N4").Value = 0.340781
$ws.Range("O4").Value = 0.002115208419356398
$ws.Range("P4").Value = 0.002115208419356398
$ws.Range("Q4").Value = 1.185263959125555
$ws.Range("R4").Value = 10.66737563213
$ws.Range("S4").Value = 0.002054098351468434
$ws.Range("T4").Value = 0.002054098351468434

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Ryk"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.43424333333333
$ws.Range("H5").Value = 31.30273
$ws.Range("I5").Value = 0.9711091978791583
$ws.Range("J5").Value = 0.9711091978791584
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 24.56783033333333
$ws.Range("N5").Value = 73.703491
$ws.Range("O5").Value = 0.4574734057918678
$ws.Range("P5").Value = 0.4574734057918678
$ws.Range("Q5").Value = 256.3467198700478
$ws.Range("R5").Value = 2307.12047883043
$ws.Range("S5").Value = 0.4442566321495875
$ws.Range("T5").Value = 0.4442566321495874

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Ryk"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.310422
$ws.Range("H6").Value = 0.9312659999999999
$ws.Range("I6").Value = 0.02889080212084161
$ws.Range("J6").Value = 0.02889080212084161
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 5.611633666666667
$ws.Range("N6").Value = 16.834901
$ws.Range("O6").Value = 0.1044932796553548
$ws.Range("P6").Value = 0.1044932796553548
$ws.Range("Q6").Value = 1.741974546074
$ws.Range("R6").Value = 15.677770914666
$ws.Range("S6").Value = 0.003018894665480618
$ws.Range("T6").Value = 0.003018894665480619

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Ryk"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.310422
$ws.Range("H7").Value = 0.9312659999999999
$ws.Range("I7").Value = 0.02889080212084161
$ws.Range("J7").Value = 0.02889080212084161
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.41023966666667
$ws.Range("N7").Value = 70.230719
$ws.Range("O7").Value = 0.435918106133421
$ws.Range("P7").Value = 0.435918106133421
$ws.Range("Q7").Value = 7.267053417805999
$ws.Range("R7").Value = 65.40348076025398
$ws.Range("S7").Value = 0.0125940237451927
$ws.Range("T7").Value = 0.0125940237451927

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Ryk"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.310422
$ws.Range("H8").Value = 0.9312659999999999
$ws.Range("I8").Value = 0.02889080212084161
$ws.Range("J8").Value = 0.02889080212084161
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1135936666666667
$ws.Range("N8").Value = 0.340781
$ws.Range("O8").Value = 0.002115208419356398
$ws.Range("P8").Value = 0.002115208419356398
$ws.Range("Q8").Value = 0.035261973194
$ws.Range("R8").Value = 0.317357758746
$ws.Range("S8").Value = 0.00006111006788796385
$ws.Range("T8").Value = 0.00006111006788796386

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Ryk"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.310422
$ws.Range("H9").Value = 0.9312659999999999
$ws.Range("I9").Value = 0.02889080212084161
$ws.Range("J9").Value = 0.02889080212084161
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 24.56783033333333
$ws.Range("N9").Value = 73.703491
$ws.Range("O9").Value = 0.4574734057918678
$ws.Range("P9").Value = 0.4574734057918678
$ws.Range("Q9").Value = 7.626395027734
$ws.Range("R9").Value = 68.637555249606
$ws.Range("S9").Value = 0.01321677364228033
$ws.Range("T9").Value = 0.01321677364228033
